# Update Name of Algo
# Update the E-column values (Algorithm imputation results) with the
# revised RandomForest output values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.1641
$ws.Range("E14").Value = 16.98610000000001
$ws.Range("E16").Value = 16.283
$ws.Range("E21").Value = 16.8001
$ws.Range("E23").Value = 16.13459999999998
$ws.Range("E25").Value = 17.05800000000001
$ws.Range("E26").Value = 16.1689
$ws.Range("E29").Value = 16.9772
$ws.Range("E40").Value = 16.61219999999999
$ws.Range("E53").Value = 16.79920000000001
$ws.Range("E57").Value = 16.7122
$ws.Range("E59").Value = 16.2477
$ws.Range("E65").Value = 17.00800000000001
$ws.Range("E69").Value = 17.34090000000003
$ws.Range("E79").Value = 18.04570000000002
$ws.Range("E83").Value = 16.5166
$ws.Range("E91").Value = 18.41840000000002
$ws.Range("E93").Value = 17.43740000000002
$ws.Range("E100").Value = 16.5506
